# Refresh the crypto price/volume snapshot (Price + Volume(1h) columns)
# with the latest scrape values. Cells are stored as plain text in the
# sheet, so each new value is entered with a leading quote-prefix to stop
# Excel from auto-converting the numeric-looking strings/percentages into
# true numbers, then the cell style is put back to Normal so no stray
# number-format override is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'289.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'40.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.45%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.046"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.84%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07288"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'4.280"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.46%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-8.08%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9202"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.56%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1154"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-8.68%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1732"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-7.30%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08645"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.52%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04175"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.33%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.32%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-1.18%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005842"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.27%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'1.34%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'0.3278"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.18%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.890"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.37%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1352"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.26%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2885"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'5.57%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03868"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-3.94%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-0.06%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003856"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-6.46%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.68%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003727"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02319"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-8.67%"
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'-6.87%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006639"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'233.25%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007695"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.85%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1273"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.16%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007381"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.80%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007062"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-14.92%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.2895"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-16.55%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006408"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.01953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-90.02%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.10%"
$ws.Range("E51").Style = "Normal"
